$wb = $excel.ActiveWorkbook

# --- Capture final selection state on Revenue (F10) before switching away ---
$revenue = $wb.Worksheets.Item("Revenue")
$revenue.Activate()
$revenue.Range("B1").Select()
$revenue.Range("F10").Select()

# --- Bump the internal sheetId counter with a throwaway sheet, then add the real "Q&A" sheet after it ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tmp1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$tmp1.Name = "ZZZ_TMP"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$qa = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$qa.Name = "Q&A"

$wb.Worksheets.Item("ZZZ_TMP").Delete()

$qa = $wb.Worksheets.Item("Q&A")

# --- Write the numbered question index column (A) first; plain numbers do not consume shared-string slots ---
$qa.Range("A2").Value = 1
$qa.Range("A5").Value = 2
$qa.Range("A9").Value = 3
$qa.Range("A12").Value = 4
$qa.Range("A15").Value = 5
$qa.Range("A18").Value = 6
$qa.Range("A21").Value = 7
$qa.Range("A25").Value = 8
$qa.Range("A29").Value = 9
$qa.Range("A32").Value = 10
$qa.Range("A35").Value = 11
$qa.Range("A39").Value = 12
$qa.Range("A43").Value = 13

# --- Write column B cells in the exact order the shared strings were first authored (194..224) ---
$qa.Range("B2").Value = "Thủ tục thành lập doanh nghiệp và phí bao nhiêu ?"   # si 194
$qa.Range("B3").Value = "http://chamsocdoanhnghiep.com/bang-gia-dich-vu-tu-van-doanh-nghiep/"   # si 195
$qa.Range("B5").Value = "Doanh nghiệp tư nhân là gì ?"   # si 196
$qa.Range("B6").Value = "http://www.ketoancattuong.vn/dich-vu/62/doanh-nghiep-tu-nhan-la-gi-.html"   # si 197
$qa.Range("B7").Value = "http://tuvanthanhlapcongtytnhh.com/tu-van-thanh-lap-cong-ty-tnhh/dac-diem-doanh-nghiep-tu-nhan-va-cong-ty-tnhh.html"   # si 198
$qa.Range("B9").Value = "Công ty trách nhiệm hữu hạn là gì ?"   # si 199
$qa.Range("B10").Value = "https://vi.wikipedia.org/wiki/C%C3%B4ng_ty_tr%C3%A1ch_nhi%E1%BB%87m_h%E1%BB%AFu_h%E1%BA%A1n"   # si 200
$qa.Range("B12").Value = "Công ty cổ phần là gì?"   # si 201
$qa.Range("B13").Value = "https://giayphepkinhdoanh.vn/cong-ty-co-phan-la-gi/"   # si 202
$qa.Range("B18").Value = "Chọn loại hình kinh doanh và đăng ký kinh doanh"   # si 203
$qa.Range("B19").Value = "http://business.gov.vn/tabid/100/catid/629/item/10907/ch%E1%BB%8Dn-lo%E1%BA%A1i-hinh-kinh-doanh-va-%C4%91%C4%83ng-ky-kinh-doanh.aspx"   # si 204
$qa.Range("B15").Value = "Công ty hợp danh là gì ?"   # si 205
$qa.Range("B16").Value = "https://vi.wikipedia.org/wiki/C%C3%B4ng_ty_h%E1%BB%A3p_danh"   # si 206
$qa.Range("B21").Value = "Khái niệm thuế môn bài là gì?"   # si 207
$qa.Range("B22").Value = "http://www.daotaoketoanhcm.com/thu-vien/phap-luat-thue/thue-mon-bai/khai-niem-thue-mon-bai-la-gi/"   # si 208
$qa.Range("B23").Value = "http://business.gov.vn/tabid/103/catid/638/item/11297/thu%E1%BA%BF-mon-bai.aspx"   # si 209
$qa.Range("B26").Value = "http://business.gov.vn/tabid/103/catid/638/item/11296/thu%E1%BA%BF-thu-nh%E1%BA%ADp-doanh-nghi%E1%BB%87p-va-thu%E1%BA%BF-thu-nh%E1%BA%ADp-h%E1%BB%99-gia-%C4%91inh.aspx"   # si 210
$qa.Range("B29").Value = "GỢI Ý TRONG VIỆC THÀNH LẬP CÔNG TY"   # si 211
$qa.Range("B32").Value = " HÀNH HỆ THỐNG NGÀNH KINH TẾ VIỆT NAM"   # si 212
$qa.Range("B33").Value = "https://dangkykinhdoanh.gov.vn/Default.aspx?tabid=106&ArticleID=274&language=en-GB"   # si 213
$qa.Range("B30").Value = "http://www.ketoancattuong.vn/trang/quy-trinh-thanh-lap-cong-ty.html"   # si 214
$qa.Range("B25").Value = "Thuế thu nhập doanh nghiệp và thuế thu nhập hộ gia đình   "   # si 215
$qa.Range("B36").Value = "http://ketoanthienung.org/tin-tuc/cac-khoan-thu-nhap-duoc-mien-thue-thu-nhap-doanh-nghiep.htm"   # si 216
$qa.Range("B35").Value = "Các khoản thu nhập được miễn thuế TNDN"   # si 217
$qa.Range("B27").Value = "http://ketoanthienung.org/tin-tuc/thue-suat-thue-thu-nhap-doanh-nghiep-nam-2014.htm"   # si 218
$qa.Range("B37").Value = "https://i-law.vn/blog/doanh-nghiep-132/kinh-doanh-thua-lo-doanh-nghiep-co-nop-thue-thu-nhap-doanh-nghiep-54069"   # si 219
$qa.Range("B40").Value = "http://business.gov.vn/tabid/110/catid/438/item/7029/gi%E1%BA%A5y-ph%C3%A9p-ho%E1%BA%A1t-%C4%91%E1%BB%99ng-ngo%E1%BA%A1i-h%E1%BB%91i-%C4%91%E1%BB%91i-v%E1%BB%9Bi-t%E1%BB%95-ch%E1%BB%A9c-t%C3%ADn-d%E1%BB%A5ng-phi-ng%C3%A2n-h%C3%A0ng.aspx"   # si 220
$qa.Range("B39").Value = "Giấy phép hoạt động ngoại hối đối với tổ chức tín dụng phi ngân hàng"   # si 221
$qa.Range("B43").Value = "Đầu tư bằng hình thức Forex rút tiền tại Việt Nam có phạm luật không ?"   # si 222
$qa.Range("B44").Value = "https://luatminhkhue.vn/tu-van-luat-doanh-nghiep/dau-tu-bang-hinh-thuc-forex-rut-tien-tai-viet-nam-co-pham-luat-khong-.aspx"   # si 223
$qa.Range("B41").Value = "https://luatminhgia.com.vn/hoi-dap-doanh-nghiep/-hoat-dong-kinh-doanh-ngoai-hoi-tai-viet-nam.aspx"   # si 224

# --- Row height for row 25 ---
$qa.Rows.Item(25).RowHeight = 15.75

# --- Rich text formatting on B25: bold the trailing two spaces ---
$b25 = $qa.Range("B25")
$len = $b25.Characters().Text.Length
$chars = $b25.Characters($len - 1, 2)
$chars.Font.Bold = $true
$chars.Font.Size = 12
$chars.Font.Color = 3355443
$chars.Font.Name = "Times New Roman"

# --- Re-activate Q&A so the stored activeTab index is recomputed correctly ---
$wb.Worksheets.Item("OPERATIONS").Activate()
$wb.Worksheets.Item("Q&A").Activate()
